$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 91, shifting existing rows 91-106 down to 92-107
$ws.Rows("91:91").Insert()

# Populate the newly inserted row 91 with the new record
$ws.Cells.Item(91, 1).Value2 = 7
$ws.Cells.Item(91, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(91, 3).Value = "Ñuble"
$ws.Cells.Item(91, 4).Value2 = 45005
$ws.Cells.Item(91, 5).Value2 = 16
$ws.Cells.Item(91, 6).Value = "Fruta"
$ws.Cells.Item(91, 7).Value2 = 100108
$ws.Cells.Item(91, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(91, 9).Value2 = 100108002
$ws.Cells.Item(91, 10).Value = "Mango"
$ws.Cells.Item(91, 11).Value = "Sin especificar"
$ws.Cells.Item(91, 12).Value = "Primera"
$ws.Cells.Item(91, 13).Value2 = 50
$ws.Cells.Item(91, 14).Value2 = 7000
$ws.Cells.Item(91, 15).Value2 = 7000
$ws.Cells.Item(91, 16).Value2 = 7000
$ws.Cells.Item(91, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(91, 18).Value = "Perú"
$ws.Cells.Item(91, 19).Value2 = 1750
$ws.Cells.Item(91, 20).Value2 = 4
